# Update the cryptocurrency price/volume table (GitHub Actions daily refresh).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
# For D-column cells whose new value looks like a plain number (e.g. "224.19"),
# the cell's NumberFormat is forced to Text ("@") first so Excel stores it as a
# string (matching the source data) instead of silently converting it to a
# numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '33.694.17'
$ws.Cells.Item(2, 5).Value = '  -1.06%  '
$ws.Cells.Item(3, 4).Value = '1.774.49'
$ws.Cells.Item(3, 5).Value = '  -1.00%  '
$ws.Cells.Item(4, 5).Value = '  +0.13%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '224.19'
$ws.Cells.Item(5, 5).Value = '  +0.85%  '
$ws.Cells.Item(6, 5).Value = '  -0.98%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.999'
$ws.Cells.Item(7, 5).Value = '  +0.08%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '32.02'
$ws.Cells.Item(8, 5).Value = '  +1.39%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.290'
$ws.Cells.Item(9, 5).Value = '  +2.19%  '
$ws.Cells.Item(10, 5).Value = '  -4.34%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0934'
$ws.Cells.Item(11, 5).Value = '  +1.44%  '
$ws.Cells.Item(12, 4).Value = '2.028.96'
$ws.Cells.Item(12, 5).Value = '  -0.98%  '
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '11.13'
$ws.Cells.Item(13, 5).Value = '  +3.82%  '
$ws.Cells.Item(14, 4).Value = '1.767.99'
$ws.Cells.Item(14, 5).Value = '  -1.34%  '
$ws.Cells.Item(15, 4).Value = '33.706.40'
$ws.Cells.Item(15, 5).Value = '  -0.85%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '0.609'
$ws.Cells.Item(16, 5).Value = '  -3.40%  '
$ws.Cells.Item(17, 5).Value = '  -2.67%  '
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '66.52'
$ws.Cells.Item(18, 5).Value = '  -2.32%  '
$ws.Cells.Item(19, 4).Value = '0.0₃0775'
$ws.Cells.Item(19, 5).Value = '  -0.97%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '238.33'
$ws.Cells.Item(20, 5).Value = '  -2.96%  '
$ws.Cells.Item(21, 5).Value = '  +0.05%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '10.58'
$ws.Cells.Item(22, 5).Value = '  -1.84%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '4.01'
$ws.Cells.Item(23, 5).Value = '  -1.85%  '
$ws.Cells.Item(24, 5).Value = '  -2.08%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '159.70'
$ws.Cells.Item(25, 5).Value = '  +0.83%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '16.09'
$ws.Cells.Item(26, 5).Value = '  -2.02%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.01'
$ws.Cells.Item(27, 5).Value = '  -0.03%  '
$ws.Cells.Item(28, 5).Value = '  +0.03%  '
$ws.Cells.Item(29, 5).Value = '  +0.22%  '
$ws.Cells.Item(30, 5).Value = '  +1.43%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.0511'
$ws.Cells.Item(31, 5).Value = '  -1.87%  '
$ws.Cells.Item(32, 5).Value = '  -2.88%  '
$ws.Cells.Item(33, 5).Value = '  -0.84%  '
$ws.Cells.Item(34, 5).Value = '  -1.65%  '
$ws.Cells.Item(35, 4).Value = '1.382.81'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.647'
$ws.Cells.Item(36, 5).Value = '  +0.58%  '
$ws.Cells.Item(37, 5).Value = '  -2.39%  '
$ws.Cells.Item(38, 5).Value = '  -1.40%  '
$ws.Cells.Item(39, 5).Value = '  +5.14%  '
$ws.Cells.Item(40, 5).Value = '  +0.68%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '78.15'
$ws.Cells.Item(41, 5).Value = '  -2.14%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.906'
$ws.Cells.Item(42, 5).Value = '  -4.21%  '
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '13.53'
$ws.Cells.Item(43, 5).Value = '  +13.83%  '
$ws.Cells.Item(44, 2).Value = 'MXToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.65'
$ws.Cells.Item(44, 5).Value = '  -2.79%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.07'
$ws.Cells.Item(45, 5).Value = '  +3.79%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.0499'
$ws.Cells.Item(46, 5).Value = '  +0.91%  '
$ws.Cells.Item(47, 5).Value = '  +12.98%  '
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '107.15'
$ws.Cells.Item(48, 5).Value = '  +1.55%  '
$ws.Cells.Item(49, 5).Value = '  -1.81%  '
$ws.Cells.Item(50, 4).Value = '1.930.04'
$ws.Cells.Item(50, 5).Value = '  -0.79%  '
$ws.Cells.Item(51, 5).Value = '  +0.17%  '
